$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Variables")
$ws2 = $wb.Worksheets.Item("Categories")

# --- Sheet "Variables": insert 4 new rows at row 6 for the new sport variables ---
$ws1.Rows.Item(6).Resize(4).Insert()

# Clear the inherited formatting on the B/C columns of the new rows - only the
# D column (valueType) should keep the red "customFormat" style copied down
# from the insert operation.
$ws1.Range("B6:C9").ClearFormats()

# Fill in the new rows. The exact order of data entry matters for the shared
# string table layout: all of column B first (top-to-bottom), then column C
# in the order 7,8,9,6, then column D.
$ws1.Range("B6").Value = "sportwi1"
$ws1.Range("B7").Value = "sportwi3"
$ws1.Range("B8").Value = "sportso1"
$ws1.Range("B9").Value = "sportso3"

$ws1.Range("C7").Value = "Which of the following items best describes your sports activity in winter time?"
$ws1.Range("C8").Value = "How often do you do sport in summer time?"
$ws1.Range("C9").Value = "Which of the following items best describes your sports activity in summer time?"
$ws1.Range("C6").Value = "How often do you do sport in winter time?"

$ws1.Range("D6").Value = "integer"
$ws1.Range("D7").Value = "integer"
$ws1.Range("D8").Value = "integer"
$ws1.Range("D9").Value = "integer"

# --- View-state updates ---
# Categories sheet: keep its existing selection (A22:A30), just nudge the
# scroll position; it should no longer be the active tab once we're done.
$ws2.Activate()
$ws2.Range("A22:A30").Select()
$excel.ActiveWindow.ScrollRow = 73

# Variables sheet: becomes the active tab, with B9 selected (cursor resting
# on the newly added sportso3 row/name cell).
$ws1.Activate()
$ws1.Range("B9").Select()
